# Apply the commit "xoa cac phan tong cua ti le chiet khau":
# 1) On "Đơn sale chính" set M3 (tỉ lệ chiết khấu sale chính, tổng) from 0.1 to 0.
# 2) On "Lương" sheet, remove the "... tại HỆ THỐNG" detail rows (rows 4-10) and the
#    resulting "Tổng lương tại HỆ THỐNG" summary row, then refresh the recomputed
#    totals that shifted as a consequence.

$wb = $excel.ActiveWorkbook

$wsSale = $wb.Worksheets.Item("Đơn sale chính")
$wsSale.Range("M3").Value = 0

$wsLuong = $wb.Worksheets.Item("Lương")

# Remove the 7 "... tại HỆ THỐNG" rows (rows 4 through 10).
$wsLuong.Range("A4:B10").EntireRow.Delete() | Out-Null

# After that deletion, the old "Tổng lương tại HỆ THỐNG" row has shifted up to row 28;
# it is no longer meaningful once HỆ THỐNG has no detail rows, so remove it too.
$wsLuong.Range("A28:B28").EntireRow.Delete() | Out-Null

# Refresh values that changed because of recomputation (ngày công, phụ cấp, lương cơ
# bản per cơ sở and the dependent "Tổng lương" rows).
$wsLuong.Range("B2").Value = 25
$wsLuong.Range("B3").Value = 875000
$wsLuong.Range("B4").Value = 2678571.428571429
$wsLuong.Range("B12").Value = 1785714.285714286
$wsLuong.Range("B20").Value = 2678571.428571429
$wsLuong.Range("B28").Value = 13571.42857142864
$wsLuong.Range("B29").Value = 1785714.285714286
$wsLuong.Range("B30").Value = 2678571.428571429
$wsLuong.Range("B31").Value = 4477857.142857143
